# Update Thresholds and Results
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9: text correction "The sentence is question" -> "The sentence is a question"
$ws.Range("C9").Value = "The sentence is a question"

# Row 2 results (E,G,H)
$ws.Range("E2").Value = 0.97910341998718586
$ws.Range("G2").Value = 0.98717948717948723
$ws.Range("H2").Value = 0.9870129870129869

# Row 6 results (E,G,H)
$ws.Range("E6").Value = 0.86051769794613375
$ws.Range("G6").Value = 0.91056910569105698
$ws.Range("H6").Value = 0.88235294117647056

# Row 7 results (E,G,H)
$ws.Range("E7").Value = 0.92088235736679336
$ws.Range("G7").Value = 0.9285714285714286
$ws.Range("H7").Value = 0.92307692307692302

# Row 9: threshold + results (D,E,G,H)
$ws.Range("D9").Value = 0.98000000000000009
$ws.Range("E9").Value = 0.86602540378443871
$ws.Range("G9").Value = 0.9
$ws.Range("H9").Value = 0.88888888888888895

# Column E picked up a best-fit width after the edits
$ws.Columns("E").AutoFit()

# Selection settles back on A1
$ws.Range("A1").Select()
